$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = "Eppendorf96"
$ws.Range("B8").Value = 77891
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Before First Read"

$ws.Range("E8").Select()
